# Default_locations.xlsx update:
#  - tweak the "Archive folder" help text (drop "and possibly finalized analyses")
#  - insert a new "Archive folder (searches)" row right after it
#  - extend formatting (blank styled rows) a bit further down the sheet
#  - leave the sheet selection on the (now shifted) "Server share" row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Reword the help text for the existing "Archive folder" row (row 5).
$ws.Range("C5").Value = "Usually a server location where raw MS data may be safely stored and backed up."

# 2. Insert a brand-new row right below it for "Archive folder (searches)".
#    Inserting at row 6 pushes "Server share" / "Fasta files" down to rows 7-8,
#    and the new row inherits the B/C column formatting from the row above.
$ws.Rows(6).Insert()
$ws.Range("A6").Value = "Archive folder (searches)"
$ws.Range("B6").Value = "B:/group/mspecgrp/Archive/Groups_temp"
$ws.Range("C6").Value = "Other server location where search results may be archived."

# 3. Extend the table's visual formatting a few rows further down, leaving a
#    gap row (row 9/10 stay empty/unused) then a header-style row 11 and
#    path/help-style blank rows 12-18.
$ws.Range("A1:C1").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B2:C2").Copy()
$ws.Range("B12:C12").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B12:C12").Copy()
$ws.Range("B13:C18").PasteSpecial(-4122)   # xlPasteFormats

# 4. Leave the active selection on the entire "Server share" row (now row 7).
$ws.Range("A7:XFD7").Select()
